$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, copying the style/format of the adjacent header cell (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the Save column values for each data row
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
